{"js": "// Replace each old cell/date value with its new value, mirroring the\n// canonical diff (one literal replacement per paragraph/table cell).\nconst body = context.document.body;\nconst replacements = [\n  [\"2025-02-24 Monday\", \"2025-02-25 Tuesday\"],\n  [\"19\u00f78=2, 3\", \"49\u00f75=9, 4\"],\n  [\"45\u00f78=5, 5\", \"13\u00f78=1, 5\"],\n  [\"85\u00f76=14, 1\", \"90\u00f77=12, 6\"],\n  [\"36\u00f75=7, 1\", \"95\u00f78=11, 7\"],\n  [\"66\u00f79=7, 3\", \"87\u00f73=29, 0\"],\n  [\"87\u00f78=10, 7\", \"81\u00f73=27, 0\"],\n  [\"22\u00f72=11, 0\", \"64\u00f75=12, 4\"],\n  [\"27\u00f78=3, 3\", \"68\u00f79=7, 5\"],\n  [\"28\u00f72=14, 0\", \"38\u00f76=6, 2\"],\n  [\"21\u00f78=2, 5\", \"23\u00f75=4, 3\"],\n  [\"59\u00f72=29, 1\", \"98\u00f72=49, 0\"],\n  [\"34\u00f75=6, 4\", \"93\u00f74=23, 1\"],\n  [\"23\u00f77=3, 2\", \"36\u00f72=18, 0\"],\n  [\"22\u00f77=3, 1\", \"89\u00f77=12, 5\"],\n  [\"43\u00f79=4, 7\", \"27\u00f79=3, 0\"],\n  [\"99\u00f75=19, 4\", \"83\u00f78=10, 3\"],\n  [\"55\u00f76=9, 1\", \"44\u00f72=22, 0\"],\n  [\"30\u00f79=3, 3\", \"29\u00f76=4, 5\"],\n  [\"25\u00f79=2, 7\", \"65\u00f73=21, 2\"],\n  [\"38\u00f79=4, 2\", \"10\u00f72=5, 0\"],\n  [\"35\u00f73=11, 2\", \"26\u00f75=5, 1\"],\n  [\"67\u00f79=7, 4\", \"10\u00f73=3, 1\"],\n  [\"64\u00f73=21, 1\", \"22\u00f73=7, 1\"],\n  [\"66\u00f74=16, 2\", \"54\u00f79=6, 0\"],\n  [\"98\u00f75=19, 3\", \"72\u00f72=36, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  // Each old string is unique in the document, so the first (only) hit is the target.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace each old cell/date value with its new value using Word's Find & Replace,\n# mirroring the canonical diff (one literal replacement per paragraph/table cell).\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-02-24 Monday', '2025-02-25 Tuesday'),\n    @('19\u00f78=2, 3', '49\u00f75=9, 4'),\n    @('45\u00f78=5, 5', '13\u00f78=1, 5'),\n    @('85\u00f76=14, 1', '90\u00f77=12, 6'),\n    @('36\u00f75=7, 1', '95\u00f78=11, 7'),\n    @('66\u00f79=7, 3', '87\u00f73=29, 0'),\n    @('87\u00f78=10, 7', '81\u00f73=27, 0'),\n    @('22\u00f72=11, 0', '64\u00f75=12, 4'),\n    @('27\u00f78=3, 3', '68\u00f79=7, 5'),\n    @('28\u00f72=14, 0', '38\u00f76=6, 2'),\n    @('21\u00f78=2, 5', '23\u00f75=4, 3'),\n    @('59\u00f72=29, 1', '98\u00f72=49, 0'),\n    @('34\u00f75=6, 4', '93\u00f74=23, 1'),\n    @('23\u00f77=3, 2', '36\u00f72=18, 0'),\n    @('22\u00f77=3, 1', '89\u00f77=12, 5'),\n    @('43\u00f79=4, 7', '27\u00f79=3, 0'),\n    @('99\u00f75=19, 4', '83\u00f78=10, 3'),\n    @('55\u00f76=9, 1', '44\u00f72=22, 0'),\n    @('30\u00f79=3, 3', '29\u00f76=4, 5'),\n    @('25\u00f79=2, 7', '65\u00f73=21, 2'),\n    @('38\u00f79=4, 2', '10\u00f72=5, 0'),\n    @('35\u00f73=11, 2', '26\u00f75=5, 1'),\n    @('67\u00f79=7, 4', '10\u00f73=3, 1'),\n    @('64\u00f73=21, 1', '22\u00f73=7, 1'),\n    @('66\u00f74=16, 2', '54\u00f79=6, 0'),\n    @('98\u00f75=19, 3', '72\u00f72=36, 0'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceOne)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
